$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks and clear all cell content/formatting,
# since the table layout is being rebuilt (16 cols x 2 rows -> 15 cols x 3 rows).
$ws.Hyperlinks.Delete()
$ws.Cells.Clear()

# Row 1: field labels (human-readable column headers)
$ws.Range("A1").Value = "タイトル"
$ws.Range("B1").Value = "説明"
$ws.Range("C1").Value = "利用条件"
$ws.Range("D1").Value = "ID"
$ws.Range("E1").Value = "ウェブサイトURL"
$ws.Range("F1").Value = "アイテムURL"
$ws.Range("G1").Value = "サムネイル"
$ws.Range("H1").Value = "機械可読ドキュメント"
$ws.Range("I1").Value = "帰属"
$ws.Range("J1").Value = "viewingDirection"
$ws.Range("K1").Value = "コレクション"
$ws.Range("L1").Value = "IIIFマニフェストURI"
$ws.Range("M1").Value = "ソート用項目"
$ws.Range("N1").Value = "西暦"
$ws.Range("O1").Value = "# of media"

# Row 2: RDF/metadata predicate names for each column
$ws.Range("A2").Value = "dcterms:title"
$ws.Range("B2").Value = "dcterms:description"
$ws.Range("C2").Value = "dcterms:rights"
$ws.Range("D2").Value = "bibo:identifier"
$ws.Range("E2").Value = "dcterms:isPartOf"
$ws.Range("F2").Value = "dcterms:relation"
$ws.Range("G2").Value = "foaf:thumbnail"
$ws.Range("H2").Value = "rdfs:seeAlso"
$ws.Range("I2").Value = "sc:attributionLabel"
$ws.Range("J2").Value = "sc:viewingDirection"
$ws.Range("K2").Value = "uterms:databaseLabel"
$ws.Range("L2").Value = "uterms:manifestUri"
$ws.Range("M2").Value = "uterms:sort"
$ws.Range("N2").Value = "uterms:year"
$ws.Range("O2").Value = 1

# Row 3: actual metadata values for this item
$ws.Range("A3").Value = "松乃栄"
$ws.Range("B3").Value = "「松乃栄(まつのさかえ)」は「旧幕府の姫君加州へ御輿入の図」という副題を持つ資料で、総合図書館に貴重書として所蔵されています。`r`nこの資料は、文政10(1827)年に徳川第11代将軍家斉の第21女・溶姫が加賀藩第13代藩主前田斉泰に輿入れしたときの様子を、三代歌川国貞が想像を交えて描いた錦絵です。東京大学のシンボルの一つである「赤門」は、このとき溶姫を迎えるため建立されたもので、白無垢の花嫁衣裳に身を包んだ溶姫が、豪奢な行列を従えて赤門をくぐる図は当時の華やかさを今に伝えています。もっとも、この資料は明治22年に描かれたもので、明治22(1889)年は家康が江戸へ入府した天正18(1590)年から数えて三百年に当たり、東京開市三百年祭が営まれた年であったため、溶姫の輿入れが描かれたと考えられます。東京大学にとっては、赤門の由来を伝える貴重な絵画史料と言えます。`r`n`r`n【請求記号 A00:6569】"
$ws.Range("C3").Value = "https://www.lib.u-tokyo.ac.jp/ja/library/general/reuse"
$ws.Range("D3").Value = "c416f868-754f-4fed-9974-6ba911e2c0ba"
$ws.Range("E3").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/matsunosakae/"
$ws.Range("F3").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/matsunosakae/document/c416f868-754f-4fed-9974-6ba911e2c0ba"
$ws.Range("G3").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif-img/21827/full/200,151/0/default.jpg"
$ws.Range("H3").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/21824"
$ws.Range("I3").Value = "東京大学総合図書館 General Library in the University of Tokyo, JAPAN"
$ws.Range("K3").Value = "松乃栄"
$ws.Range("L3").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif/c416f868-754f-4fed-9974-6ba911e2c0ba/manifest"
$ws.Range("O3").Value = 1

# Hyperlinks for URL-valued cells in row 3
$ws.Hyperlinks.Add($ws.Range("C3"), $ws.Range("C3").Value())
$ws.Hyperlinks.Add($ws.Range("E3"), $ws.Range("E3").Value())
$ws.Hyperlinks.Add($ws.Range("F3"), $ws.Range("F3").Value())
$ws.Hyperlinks.Add($ws.Range("G3"), $ws.Range("G3").Value())
$ws.Hyperlinks.Add($ws.Range("H3"), $ws.Range("H3").Value())
$ws.Hyperlinks.Add($ws.Range("L3"), $ws.Range("L3").Value())

Write-Host "done"
